# Update odds data for rows 2-6 (games on 2025-12-24) to reflect the latest
# Betfair Back/Lay snapshot, per the authoritative commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 1.45
$ws.Range("G2").Value = 1.48
$ws.Range("H2").Value = 7
$ws.Range("I2").Value = 7.8
$ws.Range("J2").Value = 5.1
$ws.Range("K2").Value = 5.5
$ws.Range("N2").Value = 4.4
$ws.Range("T2").Value = 1.98
$ws.Range("U2").Value = 1.86
$ws.Range("V2").Value = 1.15
$ws.Range("W2").Value = 3.05
$ws.Range("X2").Value = 23
$ws.Range("Y2").Value = 46
$ws.Range("AB2").Value = 8.4
$ws.Range("AC2").Value = 13
$ws.Range("AD2").Value = 1000
$ws.Range("AF2").Value = 9.4
$ws.Range("AG2").Value = 9.800000000000001
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 14.5
$ws.Range("AK2").Value = 28
$ws.Range("AN2").Value = 7.8

# Row 3
$ws.Range("F3").Value = 6.6
$ws.Range("G3").Value = 7.8
$ws.Range("H3").Value = 1.44
$ws.Range("I3").Value = 1.51
$ws.Range("J3").Value = 5
$ws.Range("K3").Value = 5.8
$ws.Range("L3").Value = 1.29
$ws.Range("M3").Value = 1.04
$ws.Range("Q3").Value = 1.56
$ws.Range("S3").Value = 2.44
$ws.Range("T3").Value = 1.82
$ws.Range("U3").Value = 2.1
$ws.Range("V3").Value = 2.96
$ws.Range("W3").Value = 1.15
$ws.Range("X3").Value = 990
$ws.Range("Y3").Value = 11.5
$ws.Range("Z3").Value = 11
$ws.Range("AA3").Value = 14
$ws.Range("AB3").Value = 50
$ws.Range("AC3").Value = 12.5
$ws.Range("AG3").Value = 44
$ws.Range("AI3").Value = 75
$ws.Range("AK3").Value = 260
$ws.Range("AL3").Value = 260
$ws.Range("AO3").Value = 5.8

# Row 4
$ws.Range("F4").Value = 1.43
$ws.Range("H4").Value = 6.8
$ws.Range("I4").Value = 20
$ws.Range("J4").Value = 4.5
$ws.Range("K4").Value = 6.2
$ws.Range("L4").Value = 1.33
$ws.Range("N4").Value = 4.6
$ws.Range("O4").Value = 1.23
$ws.Range("P4").Value = 2.22
$ws.Range("Q4").Value = 1.7
$ws.Range("R4").Value = 1.47
$ws.Range("S4").Value = 2.78
$ws.Range("W4").Value = 2.92

# Row 5
$ws.Range("F5").Value = 2.02
$ws.Range("G5").Value = 2.16
$ws.Range("H5").Value = 3.95
$ws.Range("I5").Value = 4.6
$ws.Range("J5").Value = 3.3
$ws.Range("Q5").Value = 2.1
$ws.Range("W5").Value = 1.86
$ws.Range("Z5").Value = 80
$ws.Range("AE5").Value = 1000
$ws.Range("AJ5").Value = 95
$ws.Range("AL5").Value = 1000

# Row 6
$ws.Range("F6").Value = 1.78
$ws.Range("G6").Value = 1.87
$ws.Range("H6").Value = 5.3
$ws.Range("J6").Value = 3.5
$ws.Range("K6").Value = 3.85
$ws.Range("L6").Value = 1.45
$ws.Range("N6").Value = 3.4
$ws.Range("O6").Value = 1.35
$ws.Range("P6").Value = 1.8
$ws.Range("R6").Value = 1.29
$ws.Range("S6").Value = 3.75
$ws.Range("T6").Value = 1.91
$ws.Range("U6").Value = 1.89
$ws.Range("W6").Value = 2.14
$ws.Range("X6").Value = 13.5
$ws.Range("Y6").Value = 30
$ws.Range("Z6").Value = 130
$ws.Range("AD6").Value = 46
$ws.Range("AG6").Value = 9
$ws.Range("AH6").Value = 44
$ws.Range("AJ6").Value = 40
$ws.Range("AK6").Value = 46
$ws.Range("AL6").Value = 150
$ws.Range("AN6").Value = 22

